$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(8).Delete()
$ws.Rows(7).Delete()
$ws.Rows(5).Delete()
